$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("_settings")

# Update the REST URL cell (B1) on _settings sheet
$wsSettings.Range("B1").Value = "http://10.81.21.140:8280/rill-analysis-web/rest/"

# Clear D1 and E1 values (they were [Time].[2011] helper values) but keep style
$wsSettings.Range("D1").Value = ""
$wsSettings.Range("E1").Value = ""

# Widen column B on _settings sheet (engine snaps to an MDW-7 pixel grid,
# so 53.1 is the closest input that rounds to the target stored width)
$wsSettings.Columns.Item(2).ColumnWidth = 53.1

# Make _settings the active/selected sheet (was _input before)
$wsSettings.Activate()

$wb.Save()
